{"js": "// Apply the 5 edits described by the commit \"tela index e css\".\n//\n// 1) \"...a taxa de pessoal...\"      -> \"...a taxa do pessoal...\"\n// 2) \"agronomia (procurar a palavra certa)\" -> \"recursos\"\n// 3) Underline the paragraph \"Guia pr\u00e1tico: o que \u00e9 o gr\u00e1fico e como se tornar day trader\"\n// 4) \"Trello\" paragraph -> \"O Trello \u00e9 uma plataforma de gerenciamento que...\"\n// 5) \"Semana do dia 28/ novembro de 2022\" -> \"Entrega:emana do dia 28/ novembro de 2022\"\n\nconst body = context.document.body;\n\n// --- Change 1 -----------------------------------------------------------\nlet r1 = body.search(\"taxa de pessoal\", { matchCase: true });\nr1.load(\"text\");\nawait context.sync();\nif (r1.items.length > 0) {\n  r1.items[0].insertText(\"taxa do pessoal\", \"Replace\");\n}\n\n// --- Change 2 -------------------------------------------------------------\nlet r2 = body.search(\"agronomia (procurar a palavra certa)\", { matchCase: true });\nr2.load(\"text\");\nawait context.sync();\nif (r2.items.length > 0) {\n  r2.items[0].insertText(\"recursos\", \"Replace\");\n}\n\n// --- Change 3 -------------------------------------------------------------\nconst paras = body.paragraphs;\nparas.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text === \"Guia pr\u00e1tico: o que \u00e9 o gr\u00e1fico e como se tornar day trader\") {\n    paras.items[i].font.underline = \"Single\";\n  }\n}\nawait context.sync();\n\n// --- Change 4 ---------------------------------------------------------------\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text === \"Trello\") {\n    const p = paras.items[i];\n    p.insertText(\"O \", \"Start\");\n    p.insertText(\" \u00e9 uma plataforma de gerenciamento que...\", \"End\");\n  }\n}\nawait context.sync();\n\n// --- Change 5 -----------------------------------------------------------\nlet r5 = body.search(\"Semana do dia 28/ novembro de 2022\", { matchCase: true });\nr5.load(\"text\");\nawait context.sync();\nif (r5.items.length > 0) {\n  r5.items[0].insertText(\"Entrega:emana do dia 28/ novembro de 2022\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Apply the 5 edits described by the commit \"tela index e css\".\n#\n# 1) \"...a taxa de pessoal...\"      -> \"...a taxa do pessoal...\"\n# 2) \"agronomia (procurar a palavra certa)\" -> \"recursos\"\n# 3) Underline the paragraph \"Guia pr\u00e1tico: o que \u00e9 o gr\u00e1fico e como se tornar day trader\"\n# 4) \"Trello\" paragraph -> \"O Trello \u00e9 uma plataforma de gerenciamento que...\"\n# 5) \"Semana do dia 28/ novembro de 2022\" -> \"Entrega:emana do dia 28/ novembro de 2022\"\n\n$d = $word.ActiveDocument\n\n# --- Change 1 -------------------------------------------------------------\n$rng1 = $d.Content\n$find1 = $rng1.Find\n$find1.Execute(\"taxa de pessoal\", $false, $false, $false, $false, $false, $true, 1, $false, \"taxa do pessoal\", 2) | Out-Null\n\n# --- Change 2 ---------------------------------------------------------------\n$rng2 = $d.Content\n$find2 = $rng2.Find\n$find2.Text = \"agronomia (\"\n$find2.Execute() | Out-Null\nif ($find2.Found) {\n  $rng2.Text = \"recursos\"\n}\n\n$rng2b = $d.Content\n$find2b = $rng2b.Find\n$find2b.Text = \"procurar a palavra certa)\"\n$find2b.Execute() | Out-Null\nif ($find2b.Found) {\n  $rng2b.Text = \"\"\n}\n\n# --- Change 3 ---------------------------------------------------------------\nforeach ($p in $d.Paragraphs) {\n  if ($p.Range.Text -like \"Guia pr\u00e1tico: o que \u00e9 o gr\u00e1fico e como se tornar day trader*\") {\n    $p.Range.Font.Underline = 1\n  }\n}\n\n# --- Change 4 -----------------------------------------------------------\n$rng4 = $d.Content\n$find4 = $rng4.Find\n$find4.Text = \"Trello\"\n$find4.Execute() | Out-Null\nif ($find4.Found) {\n  $after = $rng4.Duplicate\n  $after.Collapse(0)\n  $after.InsertAfter(\" \u00e9 uma plataforma de gerenciamento que...\")\n  $after.LanguageIDFarEast = \"pt-BR\"\n\n  $before = $rng4.Duplicate\n  $before.Collapse(1)\n  $before.InsertBefore(\"O \")\n  $before.LanguageIDFarEast = \"pt-BR\"\n}\n\n# --- Change 5 -----------------------------------------------------------\n$rng5 = $d.Content\n$find5 = $rng5.Find\n$find5.Execute(\"Semana do dia 28/ novembro de 2022\", $false, $false, $false, $false, $false, $true, 1, $false, \"Entrega:emana do dia 28/ novembro de 2022\", 2) | Out-Null\n"}
